# Auto-generated script: update cryptocurrency Price (D) and Volume(1h) (E) columns
# to reflect the latest GitHub Actions scrape of coinranking.com.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.135.27"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.38%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.822.82"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.69%  "
$ws.Range("E4").Value = "  -0.30%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.21"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.62%  "
$ws.Range("E6").Value = "  -0.31%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4631"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.80%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3632"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.56%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07305"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.59%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8713"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.38%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.13"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.55%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.878.86"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.06%  "
$ws.Range("E13").Value = "  +3.74%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.346"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.43%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.54"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.71%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.478"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.49%  "
$ws.Range("E17").Value = "  -0.41%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008657"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.83%  "
$ws.Range("E19").Value = "  -0.22%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "27.431.18"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.66%  "
$ws.Range("E21").Value = "  -2.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.219"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.73%  "
$ws.Range("E23").Value = "  -1.34%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.100.15"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.80%  "
$ws.Range("E25").Value = "  -1.22%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.870"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.89%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.19"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.36%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.079"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.15%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.098"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.45%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "116.14"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.44%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08912"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.20%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.959"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.41%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7365"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.92%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.457"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.06%  "
$ws.Range("E35").Value = "  -2.74%  "
$ws.Range("E36").Value = "  -0.22%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.467"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.51%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.071"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.08%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05250"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.66%  "
$ws.Range("E40").Value = "  -2.22%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.925"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.49%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.166"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.28%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5209"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.58%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1630"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.01%  "
$ws.Range("E45").Value = "  -2.97%  "
$ws.Range("E46").Value = "  -2.32%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.010"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.31%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.22"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.46%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "103.39"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.68%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.636"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.33%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06260"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.03%  "

Write-Output "Updated crypto price/volume cells"
